$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F44").Value = 0.4

$ws.Range("A49").Value = "TONE_VOICE"
$ws.Range("B49").Value = "The tone of the voice"
$ws.Range("C49").Value = 50
$ws.Range("D49").Value = 150
$ws.Range("E49").Value = "low_tone%%mid_tone%%high_tone"
$ws.Range("F49").Value = 100

$ws.Range("F47").Select()
